# Daily attendance processing - 2025-12-06 22:25:41
#
# Normalises the "Recorded By" column (G): when a cell lists multiple
# recorders separated by ", ", the literal entry "System" (exact case) is
# moved to the front and the remaining recorders are sorted alphabetically.
# Cells with a single recorder, or that are already in the correct order,
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function CaseSensitiveEquals($s1, $s2) {
    if ($s1.Length -ne $s2.Length) { return $false }
    $arr1 = $s1.ToCharArray()
    $arr2 = $s2.ToCharArray()
    for ($i = 0; $i -lt $arr1.Length; $i++) {
        if ([int]$arr1[$i] -ne [int]$arr2[$i]) {
            return $false
        }
    }
    return $true
}

function TransformRecordedBy($value) {
    $rawParts = $value.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $systemParts = @()
    $restParts = @()
    foreach ($p in $parts) {
        if (CaseSensitiveEquals $p "System") {
            $systemParts += $p
        } else {
            $restParts += $p
        }
    }

    $restSorted = $restParts | Sort-Object

    $resultParts = @()
    $resultParts += $systemParts
    $resultParts += $restSorted

    return ($resultParts -join ", ")
}

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Text
    if ($current -eq $null -or $current -eq "") {
        continue
    }
    if ($current.IndexOf(",") -lt 0) {
        continue
    }
    $updated = TransformRecordedBy $current
    if ($updated -ne $current) {
        $cell.Value = $updated
    }
}
